# Automatische test-sync: 2025-06-23 18:22:50
# Append the new mail-log entry (row 10) to the "Logs" sheet and bump the
# "IT / Technisch probleem" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")

$reply = "Beste klant,`nBedankt voor je bericht. Om je verder te kunnen helpen met het inlogprobleem, hebben we wat meer informatie nodig. Zou je alsjeblieft je gebruikersnaam kunnen doorgeven, zodat we het probleem verder kunnen onderzoeken?`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"

$logs.Range("A10").Value = "Inlogproblemen"
$logs.Range("B10").Value = "mailmind.test@zohomail.eu"
$logs.Range("C10").Value = "Ik kan niet meer inloggen op mijn account. Kunnen jullie helpen?"
$logs.Range("D10").Value = "IT / Technisch probleem"
$logs.Range("E10").Value = $reply
$logs.Range("F10").Value = "2025-06-23 18:22:19"
$logs.Range("G10").Value = "Ja"

# Writing the multi-line reply into E10 makes the host auto-grow the row
# (mirrors real Excel's content-driven row height); put it back to the
# sheet's standard height so row 10 matches the other data rows.
$logs.Rows.Item(10).AutoFit() | Out-Null

# Keep the conditional formatting ranges in sync with the new last row.
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D10"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G10"))

# Update the Dashboard pivot-style summary count for the category that the
# new entry belongs to ("IT / Technisch probleem": 2 -> 3).
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 3
